# fix: template delete work city field
# Removes the "Work city" column (I) from the user import template.
# Excel shifts the following column ("Employee type", J) left into I.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The cell comment that documents the "Employee type" column lives on J3.
# Grab its text now so it can be re-anchored on I3 after the column shift.
$oldComment = $ws.Range("J3").Comment
$commentText = $oldComment.Text()
$oldComment.Delete()

# Delete column I ("Work city") - remaining columns (incl. J) shift left.
$ws.Columns.Item(9).Delete()

# Re-create the comment on its new home cell, I3 ("Employee type").
$ws.Range("I3").AddComment($commentText)

# Restore the cursor/selection position recorded in the saved file.
$ws.Range("I6").Select() | Out-Null
